$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# 1) Split "Recurrent Neural Networks(RNNs)" into two runs:
#    "Recurrent Neural Networks(" + "RNNs) " (note trailing space added)
$target = $tr.Find("RNNs)")
$target.Text = "RNNs) "

# 2) Append a trailing space to "June 12, 2024"
$dateRange = $tr.Find("June 12, 2024")
$dateRange.Text = "June 12, 2024 "
